$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = New-Object 'object[]' 9
$rows[0] = @("ECs","Efna5","Ephb1","ECs",2,0.6666666666666666,0.4223226666666666,1.266968,0.1676547342089268,0.1676547342089269,3,1,2.814284333333333,8.442853,0.7633363087964763,0.7633363087964762,1.188536064411555,10.696824579704,0.1279769459632965,0.1279769459632965)
$rows[1] = @("ECs","Efna5","Ephb1","FAPs",2,0.6666666666666666,0.4223226666666666,1.266968,0.1676547342089268,0.1676547342089269,2,0.6666666666666666,0.1814093333333333,0.5442279999999999,0.04920481176963387,0.04920481176963386,0.07661327341155554,0.6895194607039998,0.008249419639038242,0.008249419639038242)
$rows[2] = @("ECs","Efna5","Ephb1","sCs",2,0.6666666666666666,0.4223226666666666,1.266968,0.1676547342089268,0.1676547342089269,3,1,0.6911273333333333,2.073382,0.1874588794338899,0.1874588794338899,0.2918787384195555,2.626908645776,0.03142836860659207,0.03142836860659208)
$rows[3] = @("FAPs","Efna5","Ephb1","ECs",3,1,1.874986333333333,5.624959,0.7443368783435028,0.7443368783435029,3,1,2.814284333333333,8.442853,0.7633363087964763,0.7633363087964762,5.276744663114111,47.490701968027,0.5681793652158212,0.5681793652158212)
$rows[4] = @("FAPs","Efna5","Ephb1","FAPs",3,1,1.874986333333333,5.624959,0.7443368783435028,0.7443368783435029,2,0.6666666666666666,0.1814093333333333,0.5442279999999999,0.04920481176963387,0.04920481176963386,0.3401400207391111,3.061260186652,0.03662495599208893,0.03662495599208893)
$rows[5] = @("FAPs","Efna5","Ephb1","sCs",3,1,1.874986333333333,5.624959,0.7443368783435028,0.7443368783435029,3,1,0.6911273333333333,2.073382,0.1874588794338899,0.1874588794338899,1.295854304593111,11.662688741338,0.1395325571355926,0.1395325571355927)
$rows[6] = @("sCs","Efna5","Ephb1","ECs",3,1,0.2216933333333333,0.66508,0.08800838744757017,0.08800838744757018,3,1,2.814284333333333,8.442853,0.7633363087964763,0.7633363087964762,0.6239080748044444,5.61517267324,0.06717999761735835,0.06717999761735835)
$rows[7] = @("sCs","Efna5","Ephb1","FAPs",3,1,0.2216933333333333,0.66508,0.08800838744757017,0.08800838744757018,2,0.6666666666666666,0.1814093333333333,0.5442279999999999,0.04920481176963387,0.04920481176963386,0.04021723980444444,0.36195515824,0.004330436138506698,0.004330436138506698)
$rows[8] = @("sCs","Efna5","Ephb1","sCs",3,1,0.2216933333333333,0.66508,0.08800838744757017,0.08800838744757018,3,1,0.6911273333333333,2.073382,0.1874588794338899,0.1874588794338899,0.1532183222844444,1.37896490056,0.01649795369170513,0.01649795369170513)

$r = 2
foreach ($row in $rows) {
    for ($i = 0; $i -lt $row.Length; $i++) {
        $col = $i + 1
        $ws.Cells.Item($r, $col).Value = $row[$i]
    }
    $r++
}
